$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-143 down to 43-144.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with its values.
$ws.Cells.Item(42, 1).Value = 6
$ws.Cells.Item(42, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(42, 3).Value = 'Metropolitana'
$ws.Cells.Item(42, 4).Value = 45272
$ws.Cells.Item(42, 5).Value = 13
$ws.Cells.Item(42, 6).Value = 'Fruta'
$ws.Cells.Item(42, 7).Value = 100101
$ws.Cells.Item(42, 8).Value = 'Berries'
$ws.Cells.Item(42, 9).Value = 100101008
$ws.Cells.Item(42, 10).Value = 'Mora'
$ws.Cells.Item(42, 11).Value = 'Sin especificar'
$ws.Cells.Item(42, 12).Value = 'Primera'
$ws.Cells.Item(42, 13).Value = 150
$ws.Cells.Item(42, 14).Value = 6000
$ws.Cells.Item(42, 15).Value = 6000
$ws.Cells.Item(42, 16).Value = 6000
$ws.Cells.Item(42, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(42, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(42, 19).Value = 3000
$ws.Cells.Item(42, 20).Value = 2
